# Auto-generated Excel COM-interop script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "68.190.65"
$ws.Cells.Item(2, 5).Value = "  +0.55%  "
$ws.Cells.Item(3, 4).Value = "3.833.97"
$ws.Cells.Item(4, 4).Value = "'0.999"
$ws.Cells.Item(4, 5).Value = "  -0.10%  "
$ws.Cells.Item(5, 4).Value = "'599.87"
$ws.Cells.Item(5, 5).Value = "  +0.38%  "
$ws.Cells.Item(6, 4).Value = "'171.40"
$ws.Cells.Item(6, 5).Value = "  +3.30%  "
$ws.Cells.Item(7, 4).Value = "3.832.70"
$ws.Cells.Item(8, 5).Value = "  +0.01%  "
$ws.Cells.Item(9, 5).Value = "  +0.46%  "
$ws.Cells.Item(10, 5).Value = "  +1.86%  "
$ws.Cells.Item(11, 5).Value = "  +2.73%  "
$ws.Cells.Item(12, 5).Value = "  +1.33%  "
$ws.Cells.Item(13, 4).Value = "'0.0000283"
$ws.Cells.Item(13, 5).Value = "  +14.75%  "
$ws.Cells.Item(14, 4).Value = "'36.91"
$ws.Cells.Item(14, 5).Value = "  +0.42%  "
$ws.Cells.Item(15, 4).Value = "4.485.15"
$ws.Cells.Item(15, 5).Value = "  -0.24%  "
$ws.Cells.Item(16, 4).Value = "3.874.02"
$ws.Cells.Item(16, 5).Value = "  +0.45%  "
$ws.Cells.Item(17, 4).Value = "68.250.46"
$ws.Cells.Item(17, 5).Value = "  +0.49%  "
$ws.Cells.Item(18, 4).Value = "'18.39"
$ws.Cells.Item(18, 5).Value = "  +1.77%  "
$ws.Cells.Item(19, 4).Value = "'7.45"
$ws.Cells.Item(19, 5).Value = "  +1.63%  "
$ws.Cells.Item(20, 5).Value = "  +0.84%  "
$ws.Cells.Item(21, 5).Value = "  -0.33%  "
$ws.Cells.Item(22, 4).Value = "'468.02"
$ws.Cells.Item(22, 5).Value = "  +1.21%  "
$ws.Cells.Item(23, 4).Value = "'0.731"
$ws.Cells.Item(23, 5).Value = "  +0.55%  "
$ws.Cells.Item(24, 4).Value = "'0.0000158"
$ws.Cells.Item(24, 5).Value = "  -2.60%  "
$ws.Cells.Item(25, 4).Value = "'83.44"
$ws.Cells.Item(25, 5).Value = "  +0.44%  "
$ws.Cells.Item(26, 5).Value = "  +1.73%  "
$ws.Cells.Item(27, 4).Value = "'12.14"
$ws.Cells.Item(27, 5).Value = "  +0.48%  "
$ws.Cells.Item(28, 4).Value = "'10.47"
$ws.Cells.Item(28, 5).Value = "  +5.32%  "
$ws.Cells.Item(29, 5).Value = "  +0.05%  "
$ws.Cells.Item(30, 5).Value = "  -0.22%  "
$ws.Cells.Item(31, 4).Value = "3.988.81"
$ws.Cells.Item(31, 5).Value = "  -0.24%  "
$ws.Cells.Item(32, 4).Value = "'7.75"
$ws.Cells.Item(32, 5).Value = "  +0.48%  "
$ws.Cells.Item(33, 4).Value = "'2.30"
$ws.Cells.Item(33, 5).Value = "  -0.59%  "
$ws.Cells.Item(34, 5).Value = "  +0.23%  "
$ws.Cells.Item(35, 4).Value = "'9.38"
$ws.Cells.Item(35, 5).Value = "  +1.18%  "
$ws.Cells.Item(36, 4).Value = "3.803.68"
$ws.Cells.Item(36, 5).Value = "  -0.56%  "
$ws.Cells.Item(37, 5).Value = "  +19.67%  "
$ws.Cells.Item(38, 5).Value = "  +1.02%  "
$ws.Cells.Item(39, 4).Value = "'5.96"
$ws.Cells.Item(39, 5).Value = "  +1.51%  "
$ws.Cells.Item(40, 4).Value = "'0.140"
$ws.Cells.Item(40, 5).Value = "  +0.36%  "
$ws.Cells.Item(41, 4).Value = "'1.01"
$ws.Cells.Item(41, 5).Value = "  -0.48%  "
$ws.Cells.Item(42, 4).Value = "'0.999"
$ws.Cells.Item(42, 5).Value = "  -0.07%  "
$ws.Cells.Item(43, 4).Value = "'0.319"
$ws.Cells.Item(43, 5).Value = "  +2.79%  "
$ws.Cells.Item(44, 2).Value = "USDe"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(44, 4).Value = "'1.00"
$ws.Cells.Item(44, 5).Value = "  +0.00%  "
$ws.Cells.Item(45, 2).Value = "Stacks"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(45, 4).Value = "'1.98"
$ws.Cells.Item(45, 5).Value = "  +0.71%  "
$ws.Cells.Item(46, 4).Value = "'8.75"
$ws.Cells.Item(46, 5).Value = "  +3.24%  "
$ws.Cells.Item(47, 4).Value = "'417.92"
$ws.Cells.Item(47, 5).Value = "  -1.61%  "
$ws.Cells.Item(48, 4).Value = "'0.000292"
$ws.Cells.Item(48, 5).Value = "  +6.69%  "
$ws.Cells.Item(49, 4).Value = "'46.65"
$ws.Cells.Item(49, 5).Value = "  -1.01%  "
$ws.Cells.Item(50, 2).Value = "Monero"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(50, 4).Value = "'141.98"
$ws.Cells.Item(50, 5).Value = "  -1.19%  "
$ws.Cells.Item(51, 2).Value = "VeChain"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(51, 4).Value = "'0.0359"
$ws.Cells.Item(51, 5).Value = "  +1.67%  "
